$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1049.091
$ws.Range("I12").Value = 812
$ws.Range("J12").Value = 1184.5714
$ws.Range("K12").Value = 812
$ws.Range("L12").Value = 1184.5714
$ws.Range("M12").Value = -642
$ws.Range("N12").Value = -1524.5714

$ws.Range("H18").Value = 2134.2856
$ws.Range("I18").Value = 1831.6666
$ws.Range("J18").Value = 3950
$ws.Range("K18").Value = 1831.6666
$ws.Range("L18").Value = 3950
$ws.Range("M18").Value = -1547.6666
$ws.Range("N18").Value = -4518

$ws.Range("H51").Value = 7347.5713
$ws.Range("J51").Value = 6123
$ws.Range("L51").Value = 6123
$ws.Range("N51").Value = -7091

$ws.Range("H61").Value = 699.5
$ws.Range("I61").Value = 699.5
$ws.Range("K61").Value = 2098.5
$ws.Range("M61").Value = -1926.5

$ws.Range("H64").Value = 6484.25
$ws.Range("I64").Value = 5215.6
$ws.Range("J64").Value = 8598.666999999999
$ws.Range("K64").Value = 5215.6
$ws.Range("L64").Value = 8598.666999999999
$ws.Range("M64").Value = -4967.6
$ws.Range("N64").Value = -9094.666999999999

$ws.Range("H67").Value = 6484.25
$ws.Range("I67").Value = 5215.6
$ws.Range("J67").Value = 8598.666999999999
$ws.Range("K67").Value = 5215.6
$ws.Range("L67").Value = 8598.666999999999
$ws.Range("M67").Value = -4357.6
$ws.Range("N67").Value = -10314.667

$ws.Range("H69").Value = 13732.083
$ws.Range("I69").Value = 5670.3335
$ws.Range("K69").Value = 17011.0005
$ws.Range("M69").Value = -16137.0005

$ws.Range("H72").Value = 13732.083
$ws.Range("I72").Value = 5670.3335
$ws.Range("K72").Value = 51033.0015
$ws.Range("M72").Value = -46665.0015

$ws.Range("M86").ClearContents()
$ws.Range("H86").Value = 13290
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 13290
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 13290
$ws.Range("N86").Value = -15536

$ws.Range("M89").ClearContents()
$ws.Range("H89").Value = 13290
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 13290
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 66450
$ws.Range("N89").Value = -77682

$ws.Range("H98").Value = 10417357
$ws.Range("I98").Value = 15625723
$ws.Range("K98").Value = 15625723
$ws.Range("M98").Value = -15624225

$ws.Range("H99").Value = 4237.4287
$ws.Range("I99").Value = 990.75
$ws.Range("J99").Value = 8566.333000000001
$ws.Range("K99").Value = 2972.25
$ws.Range("L99").Value = 25698.999
$ws.Range("M99").Value = -1474.25
$ws.Range("N99").Value = -28694.999

$ws.Range("H101").Value = 565.41174
$ws.Range("J101").Value = 376.26666
$ws.Range("L101").Value = 1128.79998
$ws.Range("N101").Value = -4372.79998

$ws.Range("H107").Value = 1843
$ws.Range("I107").Value = 826.1111
$ws.Range("J107").Value = 3673.4
$ws.Range("K107").Value = 826.1111
$ws.Range("L107").Value = 3673.4
$ws.Range("M107").Value = 1093.8889
$ws.Range("N107").Value = -7513.4

$ws.Range("H112").Value = 3370.7837
$ws.Range("J112").Value = 3370.7837
$ws.Range("L112").Value = 10112.3511
$ws.Range("N112").Value = -12328.3511

$ws.Range("H116").Value = 14890.214
$ws.Range("I116").Value = 6495.8887
$ws.Range("K116").Value = 6495.8887
$ws.Range("M116").Value = -3053.8887

$ws.Range("H122").Value = 10417357
$ws.Range("I122").Value = 15625723
$ws.Range("K122").Value = 46877169
$ws.Range("M122").Value = -46874719

$ws.Range("H132").Value = 1768.9615
$ws.Range("I132").Value = 1790.762
$ws.Range("J132").Value = 1677.4
$ws.Range("K132").Value = 5372.286
$ws.Range("L132").Value = 5032.200000000001
$ws.Range("M132").Value = -2842.286
$ws.Range("N132").Value = -10092.2

$ws.Range("H134").Value = 108514
$ws.Range("J134").Value = 108514
$ws.Range("L134").Value = 108514
$ws.Range("N134").Value = -118654

$ws.Range("H137").Value = 2554.9312
$ws.Range("I137").Value = 2320.4167
$ws.Range("J137").Value = 3680.6
$ws.Range("K137").Value = 6961.250100000001
$ws.Range("L137").Value = 11041.8
$ws.Range("M137").Value = -4411.250100000001
$ws.Range("N137").Value = -16141.8

$ws.Range("H138").Value = 3552.7546
$ws.Range("I138").Value = 1968.1724
$ws.Range("J138").Value = 5467.4585
$ws.Range("K138").Value = 5904.5172
$ws.Range("L138").Value = 16402.3755
$ws.Range("M138").Value = -764.5172000000002
$ws.Range("N138").Value = -26682.3755

$ws.Range("H141").Value = 3857.5715
$ws.Range("I141").Value = 3194.8948
$ws.Range("J141").Value = 10153
$ws.Range("K141").Value = 9584.6844
$ws.Range("L141").Value = 30459
$ws.Range("M141").Value = -4404.6844
$ws.Range("N141").Value = -40819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 496.75
$ws.Range("I2").Value = 496.75
$ws.Range("K2").Value = 496.75
$ws.Range("M2").Value = -383.75

$ws.Range("H32").Value = 7900.9297
$ws.Range("I32").Value = 5440.2544
$ws.Range("K32").Value = 5440.2544
$ws.Range("M32").Value = -5153.2544

$ws.Range("H43").Value = 33333
$ws.Range("J43").Value = 33333
$ws.Range("L43").Value = 33333
$ws.Range("N43").Value = -33959

$ws.Range("H61").Value = 3080235.8
$ws.Range("I61").Value = 3336454
$ws.Range("K61").Value = 3336454
$ws.Range("M61").Value = -3336242

$ws.Range("H74").Value = 2733.75
$ws.Range("I74").Value = 2654.1667
$ws.Range("J74").Value = 3450
$ws.Range("K74").Value = 2654.1667
$ws.Range("L74").Value = 3450
$ws.Range("M74").Value = -1780.1667
$ws.Range("N74").Value = -5198

$ws.Range("H77").Value = 2733.75
$ws.Range("I77").Value = 2654.1667
$ws.Range("J77").Value = 3450
$ws.Range("K77").Value = 13270.8335
$ws.Range("L77").Value = 17250
$ws.Range("M77").Value = -8902.833500000001
$ws.Range("N77").Value = -25986

$ws.Range("H102").Value = 4510.533

$ws.Range("H110").Value = 5984.421
$ws.Range("I110").Value = 6096.769
$ws.Range("K110").Value = 6096.769
$ws.Range("M110").Value = -4051.769

$ws.Range("H116").Value = 496.75
$ws.Range("I116").Value = 496.75
$ws.Range("K116").Value = 496.75
$ws.Range("M116").Value = 1797.25

$ws.Range("H122").Value = 3171.1052
$ws.Range("I122").Value = 2847.2778
$ws.Range("K122").Value = 8541.8334
$ws.Range("M122").Value = -6091.8334

$ws.Range("H132").Value = 3067.1428
$ws.Range("I132").Value = 3132.8438
$ws.Range("J132").Value = 2366.3333
$ws.Range("K132").Value = 9398.5314
$ws.Range("L132").Value = 7098.999899999999
$ws.Range("M132").Value = -6868.5314
$ws.Range("N132").Value = -12158.9999

$ws.Range("H136").Value = 3080235.8
$ws.Range("I136").Value = 3336454
$ws.Range("K136").Value = 10009362
$ws.Range("M136").Value = -10006812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 496.75
$ws.Range("I3").Value = 496.75
$ws.Range("K3").Value = 496.75
$ws.Range("M3").Value = -382.75

$ws.Range("H20").Value = 5310.9688
$ws.Range("I20").Value = 8556.714
$ws.Range("K20").Value = 8556.714
$ws.Range("M20").Value = -8309.714

$ws.Range("H53").Value = 95584.75
$ws.Range("J53").Value = 95584.75
$ws.Range("L53").Value = 95584.75
$ws.Range("N53").Value = -96732.75

$ws.Range("H86").Value = 3473.2
$ws.Range("I86").Value = 2407.238
$ws.Range("J86").Value = 5960.4443
$ws.Range("K86").Value = 2407.238
$ws.Range("L86").Value = 5960.4443
$ws.Range("M86").Value = -1284.238
$ws.Range("N86").Value = -8206.444299999999

$ws.Range("H89").Value = 3473.2
$ws.Range("I89").Value = 2407.238
$ws.Range("J89").Value = 5960.4443
$ws.Range("K89").Value = 12036.19
$ws.Range("L89").Value = 29802.2215
$ws.Range("M89").Value = -6420.189999999999
$ws.Range("N89").Value = -41034.2215

$ws.Range("H94").Value = 2059.5
$ws.Range("I94").Value = 3027.8572
$ws.Range("J94").Value = 703.8
$ws.Range("K94").Value = 3027.8572
$ws.Range("L94").Value = 703.8
$ws.Range("M94").Value = -2576.8572
$ws.Range("N94").Value = -1605.8

$ws.Range("H105").Value = 718286.3
$ws.Range("I105").Value = 1272654.1
$ws.Range("J105").Value = 5527.7144
$ws.Range("K105").Value = 1272654.1
$ws.Range("L105").Value = 5527.7144
$ws.Range("M105").Value = -1270907.1
$ws.Range("N105").Value = -9021.714400000001

$ws.Range("H107").Value = 6785
$ws.Range("I107").Value = 5241.6665
$ws.Range("J107").Value = 8637
$ws.Range("K107").Value = 5241.6665
$ws.Range("L107").Value = 8637
$ws.Range("M107").Value = -3321.6665
$ws.Range("N107").Value = -12477

$ws.Range("H134").Value = 1410.2858
$ws.Range("I134").Value = 1450.6061
$ws.Range("K134").Value = 4351.8183
$ws.Range("M134").Value = -1816.8183

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21280590
$ws.Range("I31").Value = 55558896
$ws.Range("K31").Value = 55558896
$ws.Range("M31").Value = -55558601

$ws.Range("H34").Value = 21280590
$ws.Range("I34").Value = 55558896
$ws.Range("K34").Value = 55558896
$ws.Range("M34").Value = -55558694

$ws.Range("N37").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0

$ws.Range("H58").Value = 2434.3
$ws.Range("I58").Value = 2290.25
$ws.Range("J58").Value = 2722.4
$ws.Range("K58").Value = 2290.25
$ws.Range("L58").Value = 2722.4
$ws.Range("M58").Value = -2087.25
$ws.Range("N58").Value = -3128.4

$ws.Range("H99").Value = 25328.3
$ws.Range("I99").Value = 12683.857
$ws.Range("J99").Value = 54832
$ws.Range("K99").Value = 12683.857
$ws.Range("L99").Value = 54832
$ws.Range("M99").Value = -11185.857
$ws.Range("N99").Value = -57828

$ws.Range("H105").Value = 5340
$ws.Range("I105").Value = 5340
$ws.Range("K105").Value = 5340
$ws.Range("M105").Value = -3593

$ws.Range("H126").Value = 25328.3
$ws.Range("I126").Value = 12683.857
$ws.Range("J126").Value = 54832
$ws.Range("K126").Value = 38051.571
$ws.Range("L126").Value = 164496
$ws.Range("M126").Value = -35581.571
$ws.Range("N126").Value = -169436

$ws.Range("H132").Value = 1635.3489
$ws.Range("I132").Value = 1643.7273
$ws.Range("J132").Value = 1607.7
$ws.Range("K132").Value = 4931.1819
$ws.Range("L132").Value = 4823.1
$ws.Range("M132").Value = -2401.1819
$ws.Range("N132").Value = -9883.1

$ws.Range("H134").Value = 2071.3635
$ws.Range("I134").Value = 1959.3077
$ws.Range("K134").Value = 5877.9231
$ws.Range("M134").Value = -3342.9231

$ws.Range("H136").Value = 2434.3
$ws.Range("I136").Value = 2290.25
$ws.Range("J136").Value = 2722.4
$ws.Range("K136").Value = 6870.75
$ws.Range("L136").Value = 8167.200000000001
$ws.Range("M136").Value = -4320.75
$ws.Range("N136").Value = -13267.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111.77778
$ws.Range("I2").Value = 118.818184
$ws.Range("J2").Value = 100.71429
$ws.Range("K2").Value = 712.9091040000001
$ws.Range("L2").Value = 604.28574
$ws.Range("M2").Value = -599.9091040000001
$ws.Range("N2").Value = -830.28574

$ws.Range("H3").Value = 9760.143
$ws.Range("I3").Value = 5831.3335
$ws.Range("K3").Value = 17494.0005
$ws.Range("M3").Value = -17382.0005

$ws.Range("H12").Value = 1565.2273
$ws.Range("I12").Value = 1255.3
$ws.Range("J12").Value = 1823.5
$ws.Range("K12").Value = 3765.9
$ws.Range("L12").Value = 5470.5
$ws.Range("M12").Value = -3592.9
$ws.Range("N12").Value = -5816.5

$ws.Range("H32").Value = 1686218.9
$ws.Range("J32").Value = 1823462.6
$ws.Range("L32").Value = 5470387.800000001
$ws.Range("N32").Value = -5470953.800000001

$ws.Range("H107").Value = 4140564.8
$ws.Range("I107").Value = 2086.25
$ws.Range("J107").Value = 6505409.5
$ws.Range("K107").Value = 6258.75
$ws.Range("L107").Value = 19516228.5
$ws.Range("M107").Value = -4338.75
$ws.Range("N107").Value = -19520068.5

$ws.Range("H129").Value = 11910135
$ws.Range("I129").Value = 35718108
$ws.Range("J129").Value = 6149.857
$ws.Range("K129").Value = 107154324
$ws.Range("L129").Value = 18449.571
$ws.Range("M129").Value = -107149324
$ws.Range("N129").Value = -28449.571

$ws.Range("H131").Value = 3789.1538
$ws.Range("I131").Value = 2159.158
$ws.Range("J131").Value = 8213.429
$ws.Range("K131").Value = 6477.474
$ws.Range("L131").Value = 24640.287
$ws.Range("M131").Value = -1437.474
$ws.Range("N131").Value = -34720.287

$ws.Range("H133").Value = 33053.93
$ws.Range("I133").Value = 29947.727
$ws.Range("K133").Value = 89843.181
$ws.Range("M133").Value = -84783.181

$ws.Range("H138").Value = 16163.5
$ws.Range("I138").Value = 15776
$ws.Range("K138").Value = 47328
$ws.Range("M138").Value = -42188

$ws.Range("H139").Value = 7696.1816
$ws.Range("I139").Value = 3299.3333
$ws.Range("J139").Value = 12972.4
$ws.Range("K139").Value = 9897.999899999999
$ws.Range("L139").Value = 38917.2
$ws.Range("M139").Value = -4757.999899999999
$ws.Range("N139").Value = -49197.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6834.8535
$ws.Range("I70").Value = 5791.087
$ws.Range("J70").Value = 8168.5557
$ws.Range("K70").Value = 5791.087
$ws.Range("L70").Value = 8168.5557
$ws.Range("M70").Value = -5521.087
$ws.Range("N70").Value = -8708.555700000001

$ws.Range("H73").Value = 6834.8535
$ws.Range("I73").Value = 5791.087
$ws.Range("J73").Value = 8168.5557
$ws.Range("K73").Value = 5791.087
$ws.Range("L73").Value = 8168.5557
$ws.Range("M73").Value = -4855.087
$ws.Range("N73").Value = -10040.5557

$ws.Range("H97").Value = 4304.7144
$ws.Range("I97").Value = 692.86365
$ws.Range("J97").Value = 17548.166
$ws.Range("K97").Value = 692.86365
$ws.Range("L97").Value = 17548.166
$ws.Range("M97").Value = -196.86365
$ws.Range("N97").Value = -18540.166

$ws.Range("H102").Value = 1569.4642
$ws.Range("I102").Value = 1497.8846
$ws.Range("K102").Value = 1497.8846
$ws.Range("M102").Value = 124.1153999999999

$ws.Range("H122").Value = 2541852.5
$ws.Range("J122").Value = 4955
$ws.Range("L122").Value = 14865
$ws.Range("N122").Value = -19765

$ws.Range("H126").Value = 3805.52
$ws.Range("I126").Value = 2879.125
$ws.Range("K126").Value = 8637.375
$ws.Range("M126").Value = -6167.375

$ws.Range("H132").Value = 3033127.8
$ws.Range("I132").Value = 2866
$ws.Range("J132").Value = 45456790
$ws.Range("K132").Value = 8598
$ws.Range("L132").Value = 136370370
$ws.Range("M132").Value = -6068
$ws.Range("N132").Value = -136375430

$ws.Range("H139").Value = 106979.8
$ws.Range("J139").Value = 106979.8
$ws.Range("L139").Value = 106979.8
$ws.Range("N139").Value = -117259.8

$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws.Range("H141").Value = 36171.2
$ws.Range("J141").Value = 43952
$ws.Range("L141").Value = 43952
$ws.Range("N141").Value = -54312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 65764.94
$ws.Range("I20").Value = 73200.266
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 73200.266
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -72974.266
$ws.Range("N20").Value = -10452

$ws.Range("H22").Value = 4999.3335
$ws.Range("I22").Value = 4999
$ws.Range("K22").Value = 4999
$ws.Range("M22").Value = -4704

$ws.Range("H27").Value = 4999.3335
$ws.Range("I27").Value = 4999
$ws.Range("K27").Value = 4999
$ws.Range("M27").Value = -4892

$ws.Range("H68").Value = 2606367.5
$ws.Range("I68").Value = 3789625.2
$ws.Range("J68").Value = 3200.4
$ws.Range("K68").Value = 3789625.2
$ws.Range("L68").Value = 3200.4
$ws.Range("M68").Value = -3788876.2
$ws.Range("N68").Value = -4698.4

$ws.Range("H71").Value = 2606367.5
$ws.Range("I71").Value = 3789625.2
$ws.Range("J71").Value = 3200.4
$ws.Range("K71").Value = 18948126
$ws.Range("L71").Value = 16002
$ws.Range("M71").Value = -18944382
$ws.Range("N71").Value = -23490

$ws.Range("H93").Value = 2528440.8
$ws.Range("I93").Value = 1384.3
$ws.Range("J93").Value = 4634321
$ws.Range("K93").Value = 1384.3
$ws.Range("L93").Value = 4634321
$ws.Range("M93").Value = -136.3
$ws.Range("N93").Value = -4636817

$ws.Range("N95").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0

$ws.Range("H100").Value = 8077955.5
$ws.Range("I100").Value = 3501.5
$ws.Range("K100").Value = 3501.5
$ws.Range("M100").Value = -2960.5

$ws.Range("H101").Value = 47094.5
$ws.Range("J101").Value = 47094.5
$ws.Range("L101").Value = 47094.5
$ws.Range("N101").Value = -53584.5

$ws.Range("H132").Value = 3478.9534
$ws.Range("I132").Value = 2087.7144
$ws.Range("J132").Value = 6075.933
$ws.Range("K132").Value = 6263.1432
$ws.Range("L132").Value = 18227.799
$ws.Range("M132").Value = -3733.1432
$ws.Range("N132").Value = -23287.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16333.223
$ws.Range("I62").Value = 13500
$ws.Range("J62").Value = 18599.8
$ws.Range("K62").Value = 13500
$ws.Range("L62").Value = 18599.8
$ws.Range("M62").Value = -12876
$ws.Range("N62").Value = -19847.8

$ws.Range("H65").Value = 16333.223
$ws.Range("I65").Value = 13500
$ws.Range("J65").Value = 18599.8
$ws.Range("K65").Value = 67500
$ws.Range("L65").Value = 92999
$ws.Range("M65").Value = -64380
$ws.Range("N65").Value = -99239

$ws.Range("H122").Value = 1998.2858
$ws.Range("I122").Value = 1807.5278
$ws.Range("K122").Value = 5422.5834
$ws.Range("M122").Value = -2972.5834

$ws.Range("H132").Value = 2045.1765
$ws.Range("I132").Value = 1909.0769
$ws.Range("J132").Value = 2487.5
$ws.Range("K132").Value = 5727.2307
$ws.Range("L132").Value = 7462.5
$ws.Range("M132").Value = -3197.2307
$ws.Range("N132").Value = -12522.5

$ws.Range("H136").Value = 2679.62
$ws.Range("I136").Value = 2385.6667
$ws.Range("K136").Value = 7157.000100000001
$ws.Range("M136").Value = -4607.000100000001

$ws.Range("H140").Value = 55070.43
$ws.Range("J140").Value = 57915.5
$ws.Range("L140").Value = 57915.5
$ws.Range("N140").Value = -68275.5
